$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.868.80"
$ws.Range("E2").Value = "  +2.62%  "
$ws.Range("D3").Value = "1.663.56"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0621"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").Value = "1.897.96"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "1.662.72"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.550"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "249.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.89%  "
$ws.Range("D18").Value = "27.827.02"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").Value = "0.0₃0731"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.17%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +6.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0500"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "1.437.88"
$ws.Range("E33").Value = "  -7.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.56%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.932"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("E38").Value = "  -4.18%  "
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.28%  "
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").Value = "1.807.46"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E47").Value = "  +4.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.70%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0509"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.49%  "
